# Auto-generated edit script applying Tonberry_Profits market-data refresh diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5326.3335
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H121").Value = 1426.25
$ws.Range("J121").Value = 1596.5
$ws.Range("L121").Value = 4789.5
$ws.Range("N121").Value = -8283.5
$ws.Range("H137").Value = 30877.795
$ws.Range("J137").Value = 49376.523
$ws.Range("L137").Value = 148129.569
$ws.Range("N137").Value = -153229.569
$ws.Range("H138").Value = 4002.8
$ws.Range("I138").Value = 4338.4
$ws.Range("J138").Value = 3835
$ws.Range("K138").Value = 13015.2
$ws.Range("L138").Value = 11505
$ws.Range("M138").Value = -7875.199999999999
$ws.Range("N138").Value = -21785

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1060162
$ws.Range("I2").Value = 1165678.2
$ws.Range("J2").Value = 4999.5
$ws.Range("K2").Value = 1165678.2
$ws.Range("L2").Value = 4999.5
$ws.Range("M2").Value = -1165565.2
$ws.Range("N2").Value = -5225.5
$ws.Range("H32").Value = 4437.8037
$ws.Range("I32").Value = 3633.05
$ws.Range("K32").Value = 3633.05
$ws.Range("M32").Value = -3346.05
$ws.Range("H45").Value = 1070.6923
$ws.Range("I45").Value = 1001.6667
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 1001.6667
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -624.6667
$ws.Range("N45").Value = -2653
$ws.Range("H88").Value = 3179.923
$ws.Range("J88").Value = 4114.2856
$ws.Range("L88").Value = 4114.2856
$ws.Range("N88").Value = -4926.2856
$ws.Range("H91").Value = 3179.923
$ws.Range("J91").Value = 4114.2856
$ws.Range("L91").Value = 4114.2856
$ws.Range("N91").Value = -6922.2856
$ws.Range("H110").Value = 1239.8
$ws.Range("I110").Value = 1049.75
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1049.75
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 995.25
$ws.Range("N110").Value = -6090
$ws.Range("H116").Value = 1060162
$ws.Range("I116").Value = 1165678.2
$ws.Range("J116").Value = 4999.5
$ws.Range("K116").Value = 1165678.2
$ws.Range("L116").Value = 4999.5
$ws.Range("M116").Value = -1163384.2
$ws.Range("N116").Value = -9587.5
$ws.Range("H123").Value = 65994.5
$ws.Range("J123").Value = 65994.5
$ws.Range("L123").Value = 65994.5
$ws.Range("N123").Value = -75794.5
$ws.Range("H132").Value = 2007.9807
$ws.Range("J132").Value = 2424.2964
$ws.Range("L132").Value = 7272.889200000001
$ws.Range("N132").Value = -12332.8892

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1060162
$ws.Range("I3").Value = 1165678.2
$ws.Range("J3").Value = 4999.5
$ws.Range("K3").Value = 1165678.2
$ws.Range("L3").Value = 4999.5
$ws.Range("M3").Value = -1165564.2
$ws.Range("N3").Value = -5227.5
$ws.Range("H20").Value = 1602.05
$ws.Range("I20").Value = 1374.0667
$ws.Range("J20").Value = 2286
$ws.Range("K20").Value = 1374.0667
$ws.Range("L20").Value = 2286
$ws.Range("M20").Value = -1127.0667
$ws.Range("N20").Value = -2780
$ws.Range("H99").Value = 1132.5
$ws.Range("I99").Value = 1132.5
$ws.Range("K99").Value = 1132.5
$ws.Range("M99").Value = 365.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2896.6667
$ws.Range("I31").Value = 1272.6666
$ws.Range("J31").Value = 3708.6667
$ws.Range("K31").Value = 1272.6666
$ws.Range("L31").Value = 3708.6667
$ws.Range("M31").Value = -977.6666
$ws.Range("N31").Value = -4298.6667
$ws.Range("H34").Value = 2896.6667
$ws.Range("I34").Value = 1272.6666
$ws.Range("J34").Value = 3708.6667
$ws.Range("K34").Value = 1272.6666
$ws.Range("L34").Value = 3708.6667
$ws.Range("M34").Value = -1070.6666
$ws.Range("N34").Value = -4112.6667
$ws.Range("H58").Value = 1243841.2
$ws.Range("I58").Value = 1739750.9
$ws.Range("K58").Value = 1739750.9
$ws.Range("M58").Value = -1739547.9
$ws.Range("H62").Value = 3983.6667
$ws.Range("I62").Value = 4092.4285
$ws.Range("K62").Value = 4092.4285
$ws.Range("M62").Value = -3468.4285
$ws.Range("H65").Value = 3983.6667
$ws.Range("I65").Value = 4092.4285
$ws.Range("K65").Value = 20462.1425
$ws.Range("M65").Value = -17342.1425
$ws.Range("H134").Value = 2095.875
$ws.Range("I134").Value = 1380.88
$ws.Range("J134").Value = 4649.4287
$ws.Range("K134").Value = 4142.64
$ws.Range("L134").Value = 13948.2861
$ws.Range("M134").Value = -1607.64
$ws.Range("N134").Value = -19018.2861
$ws.Range("H136").Value = 1243841.2
$ws.Range("I136").Value = 1739750.9
$ws.Range("K136").Value = 5219252.699999999
$ws.Range("M136").Value = -5216702.699999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 323.4138
$ws.Range("I5").Value = 265.83334
$ws.Range("J5").Value = 364.05884
$ws.Range("K5").Value = 797.5000200000001
$ws.Range("L5").Value = 1092.17652
$ws.Range("M5").Value = -685.5000200000001
$ws.Range("N5").Value = -1316.17652
$ws.Range("H107").Value = 2619.1
$ws.Range("I107").Value = 2094.5
$ws.Range("J107").Value = 2719.024
$ws.Range("K107").Value = 6283.5
$ws.Range("L107").Value = 8157.072
$ws.Range("M107").Value = -4363.5
$ws.Range("N107").Value = -11997.072
$ws.Range("H131").Value = 8487467
$ws.Range("J131").Value = 16396.738
$ws.Range("L131").Value = 49190.21400000001
$ws.Range("N131").Value = -59270.21400000001
$ws.Range("H135").Value = 323.4138
$ws.Range("I135").Value = 265.83334
$ws.Range("J135").Value = 364.05884
$ws.Range("K135").Value = 2392.50006
$ws.Range("L135").Value = 3276.52956
$ws.Range("M135").Value = 142.4999399999997
$ws.Range("N135").Value = -8346.529559999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2459.8
$ws.Range("J102").Value = 1844.2858
$ws.Range("L102").Value = 1844.2858
$ws.Range("N102").Value = -5088.2858
$ws.Range("H107").Value = 674.3333
$ws.Range("J107").Value = 1265.75
$ws.Range("L107").Value = 1265.75
$ws.Range("N107").Value = -5105.75
$ws.Range("H132").Value = 1676390
$ws.Range("J132").Value = 4833.1763
$ws.Range("L132").Value = 14499.5289
$ws.Range("N132").Value = -19559.5289

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4994.091
$ws.Range("I40").Value = 3082.2856
$ws.Range("J40").Value = 5886.2666
$ws.Range("K40").Value = 3082.2856
$ws.Range("L40").Value = 5886.2666
$ws.Range("M40").Value = -2946.2856
$ws.Range("N40").Value = -6158.2666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16219.857
$ws.Range("J54").Value = 16219.857
$ws.Range("L54").Value = 16219.857
$ws.Range("N54").Value = -17259.857
$ws.Range("H80").Value = 64300
$ws.Range("J80").Value = 64300
$ws.Range("L80").Value = 64300
$ws.Range("N80").Value = -66296
$ws.Range("H83").Value = 64300
$ws.Range("J83").Value = 64300
$ws.Range("L83").Value = 192900
$ws.Range("N83").Value = -202884
$ws.Range("H100").Value = 285.72726
$ws.Range("I100").Value = 284.8
$ws.Range("K100").Value = 569.6
$ws.Range("M100").Value = -28.60000000000002
$ws.Range("H107").Value = 806
$ws.Range("J107").Value = 784.1429000000001
$ws.Range("L107").Value = 2352.4287
$ws.Range("N107").Value = -6192.4287
$ws.Range("H122").Value = 42626.242
$ws.Range("I122").Value = 58199.582
$ws.Range("J122").Value = 1097.3334
$ws.Range("K122").Value = 174598.746
$ws.Range("L122").Value = 3292.0002
$ws.Range("M122").Value = -172148.746
$ws.Range("N122").Value = -8192.0002
$ws.Range("H123").Value = 47408.734
$ws.Range("J123").Value = 47408.734
$ws.Range("L123").Value = 47408.734
$ws.Range("N123").Value = -57208.734
$ws.Range("H136").Value = 11576713
$ws.Range("I136").Value = 21370014
$ws.Range("K136").Value = 64110042
$ws.Range("M136").Value = -64107492
